$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference to a guaranteed-unstyled, out-of-range cell used to reset
# each edited cell back to the default style after forcing text storage
# (NumberFormat "@") so numeric-looking strings are preserved verbatim
# as text instead of being auto-coerced into Excel numbers.
$blankStyle = $ws.Range("ZZ1").Style

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.518.62"
$ws.Range("D2").Style = $blankStyle
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.01%  "
$ws.Range("E2").Style = $blankStyle

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.822.34"
$ws.Range("D3").Style = $blankStyle

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("D4").Style = $blankStyle
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.16%  "
$ws.Range("E4").Style = $blankStyle

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "315.46"
$ws.Range("D5").Style = $blankStyle
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.37%  "
$ws.Range("E5").Style = $blankStyle

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5111"
$ws.Range("D7").Style = $blankStyle
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -5.50%  "
$ws.Range("E7").Style = $blankStyle

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3950"
$ws.Range("D8").Style = $blankStyle
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -2.34%  "
$ws.Range("E8").Style = $blankStyle

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08167"
$ws.Range("D9").Style = $blankStyle
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +6.77%  "
$ws.Range("E9").Style = $blankStyle

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "41.67"
$ws.Range("D11").Style = $blankStyle
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.57%  "
$ws.Range("E11").Style = $blankStyle

# Row 12
$ws.Range("B12").NumberFormat = "@"
$ws.Range("B12").Value = "Solana"
$ws.Range("B12").Style = $blankStyle
$ws.Range("C12").NumberFormat = "@"
$ws.Range("C12").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("C12").Style = $blankStyle
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "21.10"
$ws.Range("D12").Style = $blankStyle
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.55%  "
$ws.Range("E12").Style = $blankStyle

# Row 13
$ws.Range("B13").NumberFormat = "@"
$ws.Range("B13").Value = "Polkadot"
$ws.Range("B13").Style = $blankStyle
$ws.Range("C13").NumberFormat = "@"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("C13").Style = $blankStyle
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.318"
$ws.Range("D13").Style = $blankStyle
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.17%  "
$ws.Range("E13").Style = $blankStyle

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.002"
$ws.Range("D14").Style = $blankStyle

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.516"
$ws.Range("D15").Style = $blankStyle
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -1.71%  "
$ws.Range("E15").Style = $blankStyle

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.823.69"
$ws.Range("D16").Style = $blankStyle
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.01%  "
$ws.Range("E16").Style = $blankStyle

# Row 17
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +3.64%  "
$ws.Range("E17").Style = $blankStyle

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "92.59"
$ws.Range("D18").Style = $blankStyle
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +2.94%  "
$ws.Range("E18").Style = $blankStyle

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06656"
$ws.Range("D19").Style = $blankStyle
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.78%  "
$ws.Range("E19").Style = $blankStyle

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.78"
$ws.Range("D20").Style = $blankStyle
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.36%  "
$ws.Range("E20").Style = $blankStyle

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.090"
$ws.Range("D22").Style = $blankStyle
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.26%  "
$ws.Range("E22").Style = $blankStyle

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "28.556.83"
$ws.Range("D23").Style = $blankStyle
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.13%  "
$ws.Range("E23").Style = $blankStyle

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.40"
$ws.Range("D24").Style = $blankStyle

# Row 25
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.27%  "
$ws.Range("E25").Style = $blankStyle

# Row 26
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +3.27%  "
$ws.Range("E26").Style = $blankStyle

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "156.50"
$ws.Range("D27").Style = $blankStyle
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.82%  "
$ws.Range("E27").Style = $blankStyle

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.033.75"
$ws.Range("D28").Style = $blankStyle
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -0.02%  "
$ws.Range("E28").Style = $blankStyle

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.403"
$ws.Range("D29").Style = $blankStyle
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -2.00%  "
$ws.Range("E29").Style = $blankStyle

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "125.97"
$ws.Range("D30").Style = $blankStyle
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +1.57%  "
$ws.Range("E30").Style = $blankStyle

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.113"
$ws.Range("D31").Style = $blankStyle
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -1.03%  "
$ws.Range("E31").Style = $blankStyle

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1091"
$ws.Range("D32").Style = $blankStyle
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -1.36%  "
$ws.Range("E32").Style = $blankStyle

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.763"
$ws.Range("D33").Style = $blankStyle
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +1.52%  "
$ws.Range("E33").Style = $blankStyle

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.661"
$ws.Range("D34").Style = $blankStyle
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +0.60%  "
$ws.Range("E34").Style = $blankStyle

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.07058"
$ws.Range("D35").Style = $blankStyle
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -5.06%  "
$ws.Range("E35").Style = $blankStyle

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.2229"
$ws.Range("D36").Style = $blankStyle
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -0.30%  "
$ws.Range("E36").Style = $blankStyle

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.286"
$ws.Range("D37").Style = $blankStyle
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +1.35%  "
$ws.Range("E37").Style = $blankStyle

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.838"
$ws.Range("D39").Style = $blankStyle
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -0.26%  "
$ws.Range("E39").Style = $blankStyle

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.6324"
$ws.Range("D40").Style = $blankStyle
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +0.48%  "
$ws.Range("E40").Style = $blankStyle

# Row 41
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -0.70%  "
$ws.Range("E41").Style = $blankStyle

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.180"
$ws.Range("D42").Style = $blankStyle
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.49%  "
$ws.Range("E42").Style = $blankStyle

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.000"
$ws.Range("D43").Style = $blankStyle
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +0.07%  "
$ws.Range("E43").Style = $blankStyle

# Row 44
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +0.06%  "
$ws.Range("E44").Style = $blankStyle

# Row 45
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.49%  "
$ws.Range("E45").Style = $blankStyle

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5928"
$ws.Range("D46").Style = $blankStyle

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.730"
$ws.Range("D47").Style = $blankStyle
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +0.80%  "
$ws.Range("E47").Style = $blankStyle

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "125.21"
$ws.Range("D48").Style = $blankStyle
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.15%  "
$ws.Range("E48").Style = $blankStyle

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.991"
$ws.Range("D49").Style = $blankStyle
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.73%  "
$ws.Range("E49").Style = $blankStyle

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.189"
$ws.Range("D50").Style = $blankStyle
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -0.70%  "
$ws.Range("E50").Style = $blankStyle

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06903"
$ws.Range("D51").Style = $blankStyle
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +0.27%  "
$ws.Range("E51").Style = $blankStyle
